$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 58272.332
$ws.Range("I64").Value = 252175
$ws.Range("J64").Value = 2871.5715
$ws.Range("K64").Value = 252175
$ws.Range("L64").Value = 2871.5715
$ws.Range("M64").Value = -251927
$ws.Range("N64").Value = -3367.5715

$ws.Range("H67").Value = 58272.332
$ws.Range("I67").Value = 252175
$ws.Range("J67").Value = 2871.5715
$ws.Range("K67").Value = 252175
$ws.Range("L67").Value = 2871.5715
$ws.Range("M67").Value = -251317
$ws.Range("N67").Value = -4587.5715

$ws.Range("H99").Value = 1090.0667
$ws.Range("I99").Value = 1288.6666
$ws.Range("J99").Value = 295.66666
$ws.Range("K99").Value = 3865.9998
$ws.Range("L99").Value = 886.9999799999999
$ws.Range("M99").Value = -2367.9998
$ws.Range("N99").Value = -3882.99998

$ws.Range("H100").Value = 1314.0625
$ws.Range("I100").Value = 1072.7858
$ws.Range("J100").Value = 3003
$ws.Range("K100").Value = 1072.7858
$ws.Range("L100").Value = 3003
$ws.Range("M100").Value = -531.7858000000001
$ws.Range("N100").Value = -4085

$ws.Range("H117").Value = 35973.332
$ws.Range("J117").Value = 35973.332
$ws.Range("L117").Value = 35973.332
$ws.Range("N117").Value = -45151.332

$ws.Range("H129").Value = 2501.8333
$ws.Range("J129").Value = 2984.5
$ws.Range("L129").Value = 8953.5
$ws.Range("N129").Value = -18953.5

$ws.Range("H132").Value = 28463.639
$ws.Range("I132").Value = 4231
$ws.Range("J132").Value = 113277.875
$ws.Range("K132").Value = 12693
$ws.Range("L132").Value = 339833.625
$ws.Range("M132").Value = -10163
$ws.Range("N132").Value = -344893.625

$ws.Range("H137").Value = 4032.5227
$ws.Range("I137").Value = 1056.5714
$ws.Range("J137").Value = 6749.696
$ws.Range("K137").Value = 3169.7142
$ws.Range("L137").Value = 20249.088
$ws.Range("M137").Value = -619.7142000000003
$ws.Range("N137").Value = -25349.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27192.18
$ws.Range("I32").Value = 26426.018
$ws.Range("J32").Value = 42004.668
$ws.Range("K32").Value = 26426.018
$ws.Range("L32").Value = 42004.668
$ws.Range("M32").Value = -26139.018
$ws.Range("N32").Value = -42578.668

$ws.Range("H63").Value = 2789.16
$ws.Range("I63").Value = 2476.45
$ws.Range("J63").Value = 4040
$ws.Range("K63").Value = 2476.45
$ws.Range("L63").Value = 4040
$ws.Range("M63").Value = -1790.45
$ws.Range("N63").Value = -5412

$ws.Range("H66").Value = 2789.16
$ws.Range("I66").Value = 2476.45
$ws.Range("J66").Value = 4040
$ws.Range("K66").Value = 12382.25
$ws.Range("L66").Value = 20200
$ws.Range("M66").Value = -8950.25
$ws.Range("N66").Value = -27064

$ws.Range("H97").Value = 716.2895
$ws.Range("I97").Value = 565.6
$ws.Range("J97").Value = 1006.0769
$ws.Range("K97").Value = 565.6
$ws.Range("L97").Value = 1006.0769
$ws.Range("M97").Value = -69.60000000000002
$ws.Range("N97").Value = -1998.0769

$ws.Range("H102").Value = 18573.846
$ws.Range("I102").Value = 2882.6667
$ws.Range("J102").Value = 32023.428
$ws.Range("K102").Value = 2882.6667
$ws.Range("L102").Value = 32023.428
$ws.Range("M102").Value = -1260.6667
$ws.Range("N102").Value = -35267.428

$ws.Range("H118").Value = 41073
$ws.Range("J118").Value = 41073
$ws.Range("L118").Value = 41073
$ws.Range("N118").Value = -44387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1215.8235
$ws.Range("I94").Value = 1052.6
$ws.Range("J94").Value = 2440
$ws.Range("K94").Value = 1052.6
$ws.Range("L94").Value = 2440
$ws.Range("M94").Value = -601.5999999999999
$ws.Range("N94").Value = -3342

$ws.Range("H105").Value = 2701.5881
$ws.Range("I105").Value = 2360.5
$ws.Range("J105").Value = 3188.8572
$ws.Range("K105").Value = 2360.5
$ws.Range("L105").Value = 3188.8572
$ws.Range("M105").Value = -613.5
$ws.Range("N105").Value = -6682.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 167129.39
$ws.Range("I31").Value = 1422.3438
$ws.Range("J31").Value = 265326.16
$ws.Range("K31").Value = 1422.3438
$ws.Range("L31").Value = 265326.16
$ws.Range("M31").Value = -1127.3438
$ws.Range("N31").Value = -265916.16

$ws.Range("H34").Value = 167129.39
$ws.Range("I34").Value = 1422.3438
$ws.Range("J34").Value = 265326.16
$ws.Range("K34").Value = 1422.3438
$ws.Range("L34").Value = 265326.16
$ws.Range("M34").Value = -1220.3438
$ws.Range("N34").Value = -265730.16

$ws.Range("H62").Value = 2747.4119
$ws.Range("I62").Value = 2808.3333
$ws.Range("J62").Value = 2601.2
$ws.Range("K62").Value = 2808.3333
$ws.Range("L62").Value = 2601.2
$ws.Range("M62").Value = -2184.3333
$ws.Range("N62").Value = -3849.2

$ws.Range("H65").Value = 2747.4119
$ws.Range("I65").Value = 2808.3333
$ws.Range("J65").Value = 2601.2
$ws.Range("K65").Value = 14041.6665
$ws.Range("L65").Value = 13006
$ws.Range("M65").Value = -10921.6665
$ws.Range("N65").Value = -19246

$ws.Range("H105").Value = 3579.875
$ws.Range("I105").Value = 3579.875
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3579.875
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = -1832.875

$ws.Range("H107").Value = 595.6087
$ws.Range("I107").Value = 485.2353
$ws.Range("K107").Value = 485.2353
$ws.Range("M107").Value = 1434.7647

$ws.Range("H132").Value = 49733.965
$ws.Range("I132").Value = 2233.25
$ws.Range("J132").Value = 104020.5
$ws.Range("K132").Value = 6699.75
$ws.Range("L132").Value = 312061.5
$ws.Range("M132").Value = -4169.75
$ws.Range("N132").Value = -317121.5

$ws.Range("H134").Value = 76318.31
$ws.Range("I134").Value = 1529.3636
$ws.Range("J134").Value = 179153.12
$ws.Range("K134").Value = 4588.0908
$ws.Range("L134").Value = 537459.36
$ws.Range("M134").Value = -2053.0908
$ws.Range("N134").Value = -542529.36

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 32268512
$ws.Range("I137").Value = 5647.5
$ws.Range("J137").Value = 37048196
$ws.Range("K137").Value = 16942.5
$ws.Range("L137").Value = 111144588
$ws.Range("M137").Value = -11842.5
$ws.Range("N137").Value = -111154788

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 176134.62
$ws.Range("I80").Value = 253555.25
$ws.Range("J80").Value = 4088.7778
$ws.Range("K80").Value = 253555.25
$ws.Range("L80").Value = 4088.7778
$ws.Range("M80").Value = -252557.25
$ws.Range("N80").Value = -6084.7778

$ws.Range("H83").Value = 176134.62
$ws.Range("I83").Value = 253555.25
$ws.Range("J83").Value = 4088.7778
$ws.Range("K83").Value = 1267776.25
$ws.Range("L83").Value = 20443.889
$ws.Range("M83").Value = -1262784.25
$ws.Range("N83").Value = -30427.889

$ws.Range("H97").Value = 4886.35
$ws.Range("I97").Value = 4275.4165
$ws.Range("K97").Value = 4275.4165
$ws.Range("M97").Value = -3779.4165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1979.9
$ws.Range("I68").Value = 1474.875
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1474.875
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -725.875
$ws.Range("N68").Value = -5498

$ws.Range("H71").Value = 1979.9
$ws.Range("I71").Value = 1474.875
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 7374.375
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -3630.375
$ws.Range("N71").Value = -27488

$ws.Range("H100").Value = 1959.6
$ws.Range("I100").Value = 1959.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1959.6
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = -1418.6

$ws.Range("H104").Value = 16998
$ws.Range("J104").Value = 16998
$ws.Range("L104").Value = 16998
$ws.Range("N104").Value = -23986

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8006
$ws.Range("I15").Value = 8006
$ws.Range("K15").Value = 8006
$ws.Range("M15").Value = -7718

$ws.Range("H16").Value = 45944.25
$ws.Range("J16").Value = 45944.25
$ws.Range("L16").Value = 45944.25
$ws.Range("N16").Value = -46528.25

$ws.Range("H100").Value = 385.85
$ws.Range("I100").Value = 373.16666
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 746.33332
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -205.33332
$ws.Range("N100").Value = -2082
